$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the stat_plus/stat_minus/stat_u headers to tot_plus/tot_minus/tot_u
$ws.Range("L1").Value = "tot_plus"
$ws.Range("M1").Value = "tot_minus"
$ws.Range("N1").Value = "tot_u"

# Remove the "%syst_c" column (column O) entirely; this shifts columns
# P:Y left by one (P->O, Q->P, ..., X->W, Y->X)
$ws.Range("O1").EntireColumn.Delete()

# Update selection to match the saved view state
$ws.Range("N2").Select()
